$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextCell $ws 'D2' '61.916.83'
Set-TextCell $ws 'E2' '  -0.18%  '

Set-TextCell $ws 'D3' '2.458.05'
Set-TextCell $ws 'E3' '  -1.34%  '

Set-TextCell $ws 'E4' '  +0.00%  '

Set-TextCell $ws 'D5' '549.19'
Set-TextCell $ws 'E5' '  -0.67%  '

Set-TextCell $ws 'D6' '146.38'
Set-TextCell $ws 'E6' '  -0.41%  '

Set-TextCell $ws 'E7' '  -0.01%  '

Set-TextCell $ws 'D8' '0.586'
Set-TextCell $ws 'E8' '  -2.85%  '

Set-TextCell $ws 'D9' '2.457.53'
Set-TextCell $ws 'E9' '  -1.28%  '

Set-TextCell $ws 'E10' '  -1.65%  '

Set-TextCell $ws 'E11' '  +0.33%  '

Set-TextCell $ws 'D12' '5.42'
Set-TextCell $ws 'E12' '  +0.20%  '

Set-TextCell $ws 'D13' '0.350'
Set-TextCell $ws 'E13' '  -2.97%  '

Set-TextCell $ws 'D14' '26.03'
Set-TextCell $ws 'E14' '  -0.91%  '

Set-TextCell $ws 'D15' '2.902.50'
Set-TextCell $ws 'E15' '  -1.48%  '

Set-TextCell $ws 'D16' '0.0000168'
Set-TextCell $ws 'E16' '  +1.57%  '

Set-TextCell $ws 'D17' '61.825.26'
Set-TextCell $ws 'E17' '  -0.17%  '

Set-TextCell $ws 'D18' '2.459.93'
Set-TextCell $ws 'E18' '  -1.63%  '

Set-TextCell $ws 'D19' '10.90'
Set-TextCell $ws 'E19' '  -3.21%  '

Set-TextCell $ws 'D20' '6.99'
Set-TextCell $ws 'E20' '  -0.30%  '

Set-TextCell $ws 'D21' '4.14'
Set-TextCell $ws 'E21' '  -2.44%  '

Set-TextCell $ws 'D22' '319.91'
Set-TextCell $ws 'E22' '  -1.24%  '

Set-TextCell $ws 'D23' '0.999'

Set-TextCell $ws 'E24' '  +7.52%  '

Set-TextCell $ws 'D25' '63.93'
Set-TextCell $ws 'E25' '  -1.15%  '

Set-TextCell $ws 'D26' '0.0₃0975'
Set-TextCell $ws 'E26' '  -5.44%  '

Set-TextCell $ws 'D27' '2.583.85'
Set-TextCell $ws 'E27' '  -2.19%  '

Set-TextCell $ws 'B28' 'Binance-PegBSC-USD'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws 'D28' '1.00'
Set-TextCell $ws 'E28' '  -0.01%  '

Set-TextCell $ws 'B29' 'Fetch.AI'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D29' '1.48'
Set-TextCell $ws 'E29' '  -2.32%  '

Set-TextCell $ws 'B30' 'Aptos'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D30' '7.83'
Set-TextCell $ws 'E30' '  +1.34%  '

Set-TextCell $ws 'B31' 'Bittensor'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D31' '531.46'
Set-TextCell $ws 'E31' '  -1.94%  '

Set-TextCell $ws 'D32' '8.23'
Set-TextCell $ws 'E32' '  -3.30%  '

Set-TextCell $ws 'D33' '0.146'
Set-TextCell $ws 'E33' '  -3.27%  '

Set-TextCell $ws 'E34' '  -1.79%  '

Set-TextCell $ws 'D35' '1.62'
Set-TextCell $ws 'E35' '  +1.26%  '

Set-TextCell $ws 'D36' '5.68'
Set-TextCell $ws 'E36' '  -4.30%  '

Set-TextCell $ws 'D37' '1.00'
Set-TextCell $ws 'E37' '  +0.04%  '

Set-TextCell $ws 'D38' '4.76'
Set-TextCell $ws 'E38' '  -2.19%  '

Set-TextCell $ws 'E39' '  +0.69%  '

Set-TextCell $ws 'D40' '18.20'
Set-TextCell $ws 'E40' '  -2.45%  '

Set-TextCell $ws 'D41' '1.75'
Set-TextCell $ws 'E41' '  +2.64%  '

Set-TextCell $ws 'D42' '139.87'
Set-TextCell $ws 'E42' '  -5.03%  '

Set-TextCell $ws 'E43' '  +0.03%  '

Set-TextCell $ws 'D44' '40.38'
Set-TextCell $ws 'E44' '  -1.12%  '

Set-TextCell $ws 'D45' '2.28'
Set-TextCell $ws 'E45' '  -2.59%  '

Set-TextCell $ws 'D46' '143.23'
Set-TextCell $ws 'E46' '  -3.78%  '

Set-TextCell $ws 'D47' '3.60'
Set-TextCell $ws 'E47' '  -1.14%  '

Set-TextCell $ws 'D48' '21.44'
Set-TextCell $ws 'E48' '  -0.59%  '

Set-TextCell $ws 'D49' '0.0528'
Set-TextCell $ws 'E49' '  -2.92%  '

Set-TextCell $ws 'D50' '0.591'
Set-TextCell $ws 'E50' '  -0.47%  '

Set-TextCell $ws 'D51' '0.0933'
Set-TextCell $ws 'E51' '  -2.35%  '

